$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.932.18"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "3.696.88"
$ws.Range("E3").Value = "  +3.70%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.35"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.89"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +16.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "668.95"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.434"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +7.11%  "
$ws.Range("E9").Value = "  +5.84%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "3.694.18"
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.61"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +5.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.205"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("D15").Value = "4.385.86"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000269"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").Value = "96.618.81"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.04"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +16.87%  "
$ws.Range("D19").Value = "3.705.82"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.93"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.55"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.536"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "515.96"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.46"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000208"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +5.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.58"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +6.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.18"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.168"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +11.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.06"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.17"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +7.87%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.08"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +7.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.590"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.81"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "613.07"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.74"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +26.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.163"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +8.81%  "
$ws.Range("E42").Value = "  +7.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.96"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +8.26%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0462"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +10.11%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.18"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +9.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.430"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +25.66%  "
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.62"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.59"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +5.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.78"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +4.67%  "
